$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Give the new date cells the same date-formatted style already used by
# the rows above (e.g. C24) before filling in values, same as the existing
# rows in the "Schedule" table.
$ws.Range("C24").Copy() | Out-Null
$ws.Range("C26:C29").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 26 - D32
$ws.Range("B26").Value = "D32"
$ws.Range("C26").Value2 = 43866
$ws.Range("D26").Value = "Watched some videos about RNN"

# Day labels for rows 27 and 28 entered before their activities
$ws.Range("B27").Value = "D33"
$ws.Range("C27").Value2 = 43867

$ws.Range("B28").Value = "D34"
$ws.Range("C28").Value2 = 43868

$ws.Range("D27").Value = "Completed lessons 8.4 and 8.5"
$ws.Range("D28").Value = "Completed lessons 8.6 to 8.8"

# Row 29 - D35
$ws.Range("B29").Value = "D35"
$ws.Range("C29").Value2 = 43869
$ws.Range("D29").Value = "Completed lesson 8. Participated in study jam from 10.30 AM to midnight. Submitted solutions for binary classification and style transfer problems. Took part in style transfer quizzes. Also, participated in fun activities."

# Update the selected cell to match the final state of the workbook
$ws.Range("D29").Select()
